# Increment the "想去人数" (want-to-go count) values by 1 for a few events
# on both the "展览" and "全部类型" sheets, matching the latest scrape.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F7").Value  = 48
    $ws.Range("F11").Value = 4678
    $ws.Range("F15").Value = 23
}
